$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = 0.6867480601921443
    "C2" = 0.06186139118783274
    "E2" = 0.06414771136651787
    "F2" = 4.745424882012145
    "G2" = 0.00256535435728528
    "J2" = 0.351505183128495
    "K2" = 0.6705353803752985
    "B3" = 0.6598999962387211
    "C3" = 0.05859817479149854
    "E3" = 0.06919183598000878
    "F3" = 4.558725301135297
    "G3" = 0.002570151694943194
    "J3" = 0.3304764664003699
    "K3" = 0.6429127507293515
    "B4" = 0.6440199493475234
    "C4" = 0.05667445641604729
    "E4" = 0.07261077449990871
    "F4" = 4.445366689240643
    "G4" = 0.002573248420380762
    "J4" = 0.3175583180250072
    "K4" = 0.6265921626335569
    "B5" = 0.6377003233068024
    "C5" = 0.05591047194965881
    "E5" = 0.07408422258839664
    "F5" = 4.399485825492718
    "G5" = 0.002574548506204794
    "J5" = 0.3122915013209706
    "K5" = 0.6201016840448972
    "B6" = 0.6366600988854998
    "C6" = 0.05578481336775098
    "E6" = 0.07433371058955096
    "F6" = 4.391886083729844
    "G6" = 0.002574766692285068
    "J6" = 0.3114167678251505
    "K6" = 0.6190336093477242
    "B7" = 0.6439341073448475
    "C7" = 0.05666407248553185
    "E7" = 0.07263032215313636
    "F7" = 4.444746663102535
    "G7" = 0.002573265799225455
    "J7" = 0.3174872995153351
    "K7" = 0.6265039815280602
    "B8" = 0.6773650392745765
    "C8" = 0.06071955799819762
    "E8" = 0.0658196935226858
    "F8" = 4.680781746872668
    "G8" = 0.002566977193998751
    "J8" = 0.344255237641363
    "K8" = 0.6608779299050695
    "B9" = 0.7477473729219923
    "C9" = 0.0693137209656669
    "E9" = 0.05505039898215003
    "F9" = 5.154114576379527
    "G9" = 0.002555838119578146
    "J9" = 0.3967412620999227
    "K9" = 0.7333941401487039
    "B10" = 0.8024417697434387
    "C10" = 0.07602995096017651
    "E10" = 0.04876060072384014
    "F10" = 5.508773418216691
    "G10" = 0.002548372489387998
    "J10" = 0.4353652654044424
    "K10" = 0.7898391855967191
    "B11" = 0.8279818279525841
    "C11" = 0.07917515858270008
    "E11" = 0.04626157676127285
    "F11" = 5.671731534183152
    "G11" = 0.002545130223629368
    "J11" = 0.4529644988392647
    "K11" = 0.8162174526749482
    "B12" = 0.8377486938086065
    "C12" = 0.08037929624285312
    "E12" = 0.04536820249932738
    "F12" = 5.733681661555579
    "G12" = 0.002543924443805077
    "J12" = 0.4596341701866891
    "K12" = 0.8263079116789527
    "B13" = 0.8356409718927864
    "C13" = 0.08011937741984809
    "E13" = 0.04555823767044309
    "F13" = 5.720328723847729
    "G13" = 0.002544183153921686
    "J13" = 0.4581974878404083
    "K13" = 0.8241302205719876
    "B14" = 0.8287834384573216
    "C14" = 0.07927395980952667
    "E14" = 0.04618701208945275
    "F14" = 5.676823322525308
    "G14" = 0.002545030583460238
    "J14" = 0.4535131059235766
    "K14" = 0.8170455605281575
    "B15" = 0.8245954430691427
    "C15" = 0.07875783057610874
    "E15" = 0.04657907635425573
    "F15" = 5.650206694107624
    "G15" = 0.002545552519088536
    "J15" = 0.4506444986044471
    "K15" = 0.8127192535649499
    "B16" = 0.8007859839915739
    "C16" = 0.07582622973738751
    "E16" = 0.04893128032241734
    "F16" = 5.498157057249671
    "G16" = 0.002548587464195197
    "J16" = 0.4342157872402197
    "K16" = 0.7881294768129123
    "B17" = 0.7863489311688738
    "C17" = 0.0740509551160784
    "E17" = 0.05046766751988585
    "F17" = 5.405300742987862
    "G17" = 0.002550488623916852
    "J17" = 0.4241454746395448
    "K17" = 0.773224555508591
    "B18" = 0.778107140274102
    "C18" = 0.0730383206611549
    "E18" = 0.05138539969897415
    "F18" = 5.35204442188089
    "G18" = 0.002551596613445345
    "J18" = 0.4183559609572853
    "K18" = 0.7647175861235382
    "B19" = 0.7753272457323703
    "C19" = 0.07269690705497567
    "E19" = 0.05170194891984892
    "F19" = 5.334038636029021
    "G19" = 0.002551974252881131
    "J19" = 0.4163961566662806
    "K19" = 0.7618485761599345
    "B20" = 0.7878793570276628
    "C20" = 0.07423905962319566
    "E20" = 0.05030058698194217
    "F20" = 5.415169647904321
    "G20" = 0.002550284743520793
    "J20" = 0.4252171932722888
    "K20" = 0.7748043791324051
    "B21" = 0.8307950682997216
    "C21" = 0.07952192204425046
    "E21" = 0.04600088174604089
    "F21" = 5.689595295315144
    "G21" = 0.002544781076835906
    "J21" = 0.454888872189315
    "K21" = 0.8191237338494375
    "B22" = 0.8593993206592643
    "C22" = 0.08305111329225667
    "E22" = 0.04349983197404761
    "F22" = 5.870359438007995
    "G22" = 0.002541312258334102
    "J22" = 0.4743119772812463
    "K22" = 0.8486814593164524
    "B23" = 0.8440815852742674
    "C23" = 0.08116045408495154
    "E23" = 0.04480612133818695
    "F23" = 5.773750353198864
    "G23" = 0.002543151950706138
    "J23" = 0.4639423160952845
    "K23" = 0.832851477350232
    "B24" = 0.7871872700582117
    "C24" = 0.07415399268717238
    "E24" = 0.05037601687190829
    "F24" = 5.410707517862761
    "G24" = 0.002550376871085773
    "J24" = 0.4247326692444631
    "K24" = 0.7740899475930405
    "B25" = 0.7281861820189022
    "C25" = 0.0669189759021549
    "E25" = 0.05768222066811468
    "F25" = 5.024891585113039
    "G25" = 0.002558724757183137
    "J25" = 0.3825362410441357
    "K25" = 0.6705353803752985
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
